$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 348, shifting the existing rows 348-366 down to 350-368.
$ws.Rows("348:349").Insert()

# Fill in the data for the two newly inserted rows (348-349).
$ws.Range("A348").Value = 3
$ws.Range("B348").Value = "Femacal de La Calera"
$ws.Range("C348").Value = "Coquimbo"
$ws.Range("D348").Value = 44516
$ws.Range("E348").Value = 5
$ws.Range("F348").Value = 100114001
$ws.Range("G348").Value = "Papa"
$ws.Range("H348").Value = "Rosara"
$ws.Range("I348").Value = "1a nueva(o)"
$ws.Range("J348").Value = 510
$ws.Range("K348").Value = 9000
$ws.Range("L348").Value = 9500
$ws.Range("M348").Value = 9255
$ws.Range("N348").Value = "`$/saco 25 kilos"
$ws.Range("O348").Value = "Provincia de Quillota"
$ws.Range("P348").Value = 370
$ws.Range("Q348").Value = 25
$ws.Range("R348").Value = "Hortaliza"

$ws.Range("A349").Value = 3
$ws.Range("B349").Value = "Femacal de La Calera"
$ws.Range("C349").Value = "Coquimbo"
$ws.Range("D349").Value = 44516
$ws.Range("E349").Value = 5
$ws.Range("F349").Value = 100114001
$ws.Range("G349").Value = "Papa"
$ws.Range("H349").Value = "Rosara"
$ws.Range("I349").Value = "1a nueva(o)"
$ws.Range("J349").Value = 160
$ws.Range("K349").Value = 9000
$ws.Range("L349").Value = 9000
$ws.Range("M349").Value = 9000
$ws.Range("N349").Value = "`$/saco 25 kilos"
$ws.Range("O349").Value = "Provincia de Talca"
$ws.Range("P349").Value = 360
$ws.Range("Q349").Value = 25
$ws.Range("R349").Value = "Hortaliza"
